$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new value (all values are text strings in this sheet)
$changes = @{
    "D2"  = "307.60"
    "E2"  = "-0.25%"
    "D3"  = "41.05"
    "E3"  = "0.62%"
    "D4"  = "5.240"
    "E4"  = "2.41%"
    "D5"  = "0.07666"
    "E5"  = "0.71%"
    "D6"  = "1.640"
    "E6"  = "1.09%"
    "D7"  = "0.9153"
    "E7"  = "1.66%"
    "D8"  = "2.444"
    "E8"  = "-0.16%"
    "D9"  = "0.1247"
    "E9"  = "15.03%"
    "D10" = "0.1826"
    "E10" = "3.76%"
    "D11" = "0.09093"
    "E11" = "-1.09%"
    "D12" = "0.04266"
    "E12" = "1.92%"
    "E13" = "-0.07%"
    "D14" = "0.001261"
    "E14" = "0.74%"
    "D15" = "0.005758"
    "E15" = "-2.35%"
    "D17" = "3.347"
    "E17" = "-0.17%"
    "D18" = "4.311"
    "E18" = "1.33%"
    "D20" = "7.316"
    "E20" = "11.69%"
    "D21" = "0.1384"
    "E21" = "1.41%"
    "E22" = "1.28%"
    "D23" = "0.04074"
    "E23" = "-0.36%"
    "E24" = "3.24%"
    "D25" = "0.004289"
    "E25" = "4.93%"
    "E26" = "-2.16%"
    "D38" = "0.02472"
    "E38" = "4.43%"
    "D39" = "0.05291"
    "E39" = "2.12%"
    "D40" = "0.007843"
    "E40" = "0.88%"
    "D41" = "0.1313"
    "E41" = "1.02%"
    "D42" = "0.006881"
    "E42" = "1.42%"
    "D43" = "0.001914"
    "E43" = "-1.90%"
    "D44" = "0.007638"
    "E44" = "-10.72%"
    "D45" = "0.3061"
    "E45" = "-0.43%"
    "D46" = "0.00006729"
    "E46" = "-3.11%"
    "E47" = "0.12%"
    "D48" = "0.4390"
    "E48" = "1,267.65%"
    "D49" = "0.003107"
    "E49" = "-26.10%"
    "E50" = "0.12%"
    "D51" = "0.0002004"
    "E51" = "0.12%"
}

foreach ($addr in $changes.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $changes[$addr]
}
